$wb = $excel.ActiveWorkbook

# diff hunk @1520 -> ALC!row18
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1650
$ws.Cells.Item(18, 9).Value = 2100
$ws.Cells.Item(18, 10).Value = 300
$ws.Cells.Item(18, 11).Value = 2100
$ws.Cells.Item(18, 12).Value = 300
$ws.Cells.Item(18, 13).Value = -1816
$ws.Cells.Item(18, 14).Value = -868

# diff hunk @2157 -> ALC!row31
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(31, 8).Value = 7777
$ws.Cells.Item(31, 9).Value = 7777
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 23331
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = -23101

# diff hunk @2258 -> ALC!row33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 106.61905
$ws.Cells.Item(33, 9).Value = 107.36842
$ws.Cells.Item(33, 10).Value = 99.5
$ws.Cells.Item(33, 11).Value = 107.36842
$ws.Cells.Item(33, 12).Value = 99.5
$ws.Cells.Item(33, 13).Value = 121.63158
$ws.Cells.Item(33, 14).Value = -557.5

# diff hunk @2613 -> ALC!row40
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 3942.08
$ws.Cells.Item(40, 9).Value = 3550.182
$ws.Cells.Item(40, 10).Value = 4250
$ws.Cells.Item(40, 11).Value = 3550.182
$ws.Cells.Item(40, 12).Value = 4250
$ws.Cells.Item(40, 13).Value = -3375.182
$ws.Cells.Item(40, 14).Value = -4600

# diff hunk @7881 -> ARM!row4
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = $null

# diff hunk @9226 -> ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 193227.72
$ws.Cells.Item(32, 9).Value = 196307.48
$ws.Cells.Item(32, 10).Value = 30000
$ws.Cells.Item(32, 11).Value = 196307.48
$ws.Cells.Item(32, 12).Value = 30000
$ws.Cells.Item(32, 13).Value = -196020.48
$ws.Cells.Item(32, 14).Value = -30574

# diff hunk @9866 -> ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1613.2354
$ws.Cells.Item(45, 9).Value = 1452.0834
$ws.Cells.Item(45, 10).Value = 2000
$ws.Cells.Item(45, 11).Value = 1452.0834
$ws.Cells.Item(45, 12).Value = 2000
$ws.Cells.Item(45, 13).Value = -1075.0834
$ws.Cells.Item(45, 14).Value = -2754

# diff hunk @10347 -> ARM!row55
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(55, 8).Value = 20048
$ws.Cells.Item(55, 9).Value = 20048
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 11).Value = 20048
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 13).Value = -19733

# diff hunk @10635 -> ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 16672456
$ws.Cells.Item(61, 9).Value = 6946.9
$ws.Cells.Item(61, 10).Value = 100000000
$ws.Cells.Item(61, 11).Value = 6946.9
$ws.Cells.Item(61, 12).Value = 100000000
$ws.Cells.Item(61, 13).Value = -6734.9
$ws.Cells.Item(61, 14).Value = -100000424

# diff hunk @11251 -> ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 4295792
$ws.Cells.Item(74, 9).Value = 6187255
$ws.Cells.Item(74, 10).Value = 39999.5
$ws.Cells.Item(74, 11).Value = 6187255
$ws.Cells.Item(74, 12).Value = 39999.5
$ws.Cells.Item(74, 13).Value = -6186381
$ws.Cells.Item(74, 14).Value = -41747.5

# diff hunk @11398 -> ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 4295792
$ws.Cells.Item(77, 9).Value = 6187255
$ws.Cells.Item(77, 10).Value = 39999.5
$ws.Cells.Item(77, 11).Value = 30936275
$ws.Cells.Item(77, 12).Value = 199997.5
$ws.Cells.Item(77, 13).Value = -30931907
$ws.Cells.Item(77, 14).Value = -208733.5

# diff hunk @11545 -> ARM!row80
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(80, 8).Value = 71000
$ws.Cells.Item(80, 9).Value = 42000
$ws.Cells.Item(80, 10).Value = 100000
$ws.Cells.Item(80, 11).Value = 42000
$ws.Cells.Item(80, 12).Value = 100000
$ws.Cells.Item(80, 13).Value = -41002
$ws.Cells.Item(80, 14).Value = -101996

# diff hunk @11689 -> ARM!row83
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(83, 8).Value = 71000
$ws.Cells.Item(83, 9).Value = 42000
$ws.Cells.Item(83, 10).Value = 100000
$ws.Cells.Item(83, 11).Value = 126000
$ws.Cells.Item(83, 12).Value = 300000
$ws.Cells.Item(83, 13).Value = -121008
$ws.Cells.Item(83, 14).Value = -309984

# diff hunk @14042 -> ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 1002590
$ws.Cells.Item(132, 9).Value = 1391008.4
$ws.Cells.Item(132, 10).Value = 3800
$ws.Cells.Item(132, 11).Value = 4173025.2
$ws.Cells.Item(132, 12).Value = 11400
$ws.Cells.Item(132, 13).Value = -4170495.2
$ws.Cells.Item(132, 14).Value = -16460

# diff hunk @14241 -> ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 16672456
$ws.Cells.Item(136, 9).Value = 6946.9
$ws.Cells.Item(136, 10).Value = 100000000
$ws.Cells.Item(136, 11).Value = 20840.7
$ws.Cells.Item(136, 12).Value = 300000000
$ws.Cells.Item(136, 13).Value = -18290.7
$ws.Cells.Item(136, 14).Value = -300005100

# diff hunk @15600 -> BSM!row22
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 400
$ws.Cells.Item(22, 9).Value = 400
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 400
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -227

# diff hunk @21762 -> CRP!row7
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 18275.836
$ws.Cells.Item(7, 9).Value = 55600.723
$ws.Cells.Item(7, 10).Value = 117.78378
$ws.Cells.Item(7, 11).Value = 55600.723
$ws.Cells.Item(7, 12).Value = 117.78378
$ws.Cells.Item(7, 13).Value = -55487.723
$ws.Cells.Item(7, 14).Value = -343.78378

# diff hunk @22497 -> CRP!row22
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 93939.45
$ws.Cells.Item(22, 9).Value = 267.6
$ws.Cells.Item(22, 10).Value = 171999.33
$ws.Cells.Item(22, 11).Value = 267.6
$ws.Cells.Item(22, 12).Value = 171999.33
$ws.Cells.Item(22, 13).Value = 82.39999999999998
$ws.Cells.Item(22, 14).Value = -172699.33

# diff hunk @23458 -> CRP!row41
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value = 9481
$ws.Cells.Item(41, 9).Value = 3929.5
$ws.Cells.Item(41, 10).Value = 15032.5
$ws.Cells.Item(41, 11).Value = 3929.5
$ws.Cells.Item(41, 12).Value = 15032.5
$ws.Cells.Item(41, 13).Value = -3501.5
$ws.Cells.Item(41, 14).Value = -15888.5

# diff hunk @25176 -> CRP!row76
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(76, 8).Value = 8999
$ws.Cells.Item(76, 9).Value = 8999
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 11).Value = 8999
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 13).Value = -8684

# diff hunk @25320 -> CRP!row79
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(79, 8).Value = 8999
$ws.Cells.Item(79, 9).Value = 8999
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 11).Value = 8999
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 13).Value = -7907

# diff hunk @26291 -> CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 56570.066
$ws.Cells.Item(99, 9).Value = 51363.453
$ws.Cells.Item(99, 10).Value = 70888.25
$ws.Cells.Item(99, 11).Value = 51363.453
$ws.Cells.Item(99, 12).Value = 70888.25
$ws.Cells.Item(99, 13).Value = -49865.453
$ws.Cells.Item(99, 14).Value = -73884.25

# diff hunk @27617 -> CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 56570.066
$ws.Cells.Item(126, 9).Value = 51363.453
$ws.Cells.Item(126, 10).Value = 70888.25
$ws.Cells.Item(126, 11).Value = 154090.359
$ws.Cells.Item(126, 12).Value = 212664.75
$ws.Cells.Item(126, 13).Value = -151620.359
$ws.Cells.Item(126, 14).Value = -217604.75

# diff hunk @27917 -> CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 4048
$ws.Cells.Item(132, 9).Value = 4002.0908
$ws.Cells.Item(132, 10).Value = 4174.25
$ws.Cells.Item(132, 11).Value = 12006.2724
$ws.Cells.Item(132, 12).Value = 12522.75
$ws.Cells.Item(132, 13).Value = -9476.2724
$ws.Cells.Item(132, 14).Value = -17582.75

# diff hunk @28018 -> CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 2759
$ws.Cells.Item(134, 9).Value = 2498.4
$ws.Cells.Item(134, 10).Value = 4322.6
$ws.Cells.Item(134, 11).Value = 7495.200000000001
$ws.Cells.Item(134, 12).Value = 12967.8
$ws.Cells.Item(134, 13).Value = -4960.200000000001
$ws.Cells.Item(134, 14).Value = -18037.8

# diff hunk @30138 -> CUL!row34
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 12.6
$ws.Cells.Item(34, 9).Value = 12.6
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 37.8
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = 46.2

# diff hunk @31191 -> CUL!row55
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 420000000
$ws.Cells.Item(55, 9).Value = 420000000
$ws.Cells.Item(55, 10).Value = 0
$ws.Cells.Item(55, 11).Value = 1260000000
$ws.Cells.Item(55, 12).Value = 0
$ws.Cells.Item(55, 13).Value = -1259999823
$ws.Cells.Item(55, 14).Value = $null

# diff hunk @35050 -> CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 6724.4116
$ws.Cells.Item(131, 9).Value = 1765
$ws.Cells.Item(131, 10).Value = 7385.6665
$ws.Cells.Item(131, 11).Value = 5295
$ws.Cells.Item(131, 12).Value = 22156.9995
$ws.Cells.Item(131, 13).Value = -255
$ws.Cells.Item(131, 14).Value = -32236.9995

# diff hunk @35716 -> GSM!row2
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 5316087
$ws.Cells.Item(2, 9).Value = 5941444
$ws.Cells.Item(2, 10).Value = 549.5
$ws.Cells.Item(2, 11).Value = 5941444
$ws.Cells.Item(2, 12).Value = 549.5
$ws.Cells.Item(2, 13).Value = -5941331
$ws.Cells.Item(2, 14).Value = -775.5

# diff hunk @38126 -> GSM!row51
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(51, 8).Value = 69999
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 69999
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 69999
$ws.Cells.Item(51, 14).Value = -71017

# diff hunk @39036 -> GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 15317.083
$ws.Cells.Item(70, 9).Value = 7482.8887
$ws.Cells.Item(70, 10).Value = 38819.668
$ws.Cells.Item(70, 11).Value = 7482.8887
$ws.Cells.Item(70, 12).Value = 38819.668
$ws.Cells.Item(70, 13).Value = -7212.8887
$ws.Cells.Item(70, 14).Value = -39359.668

# diff hunk @39183 -> GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 15317.083
$ws.Cells.Item(73, 9).Value = 7482.8887
$ws.Cells.Item(73, 10).Value = 38819.668
$ws.Cells.Item(73, 11).Value = 7482.8887
$ws.Cells.Item(73, 12).Value = 38819.668
$ws.Cells.Item(73, 13).Value = -6546.8887
$ws.Cells.Item(73, 14).Value = -40691.668

# diff hunk @40822 -> GSM!row107
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 5538.9
$ws.Cells.Item(107, 9).Value = 6121.8335
$ws.Cells.Item(107, 10).Value = 292.5
$ws.Cells.Item(107, 11).Value = 6121.8335
$ws.Cells.Item(107, 12).Value = 292.5
$ws.Cells.Item(107, 13).Value = -4201.8335
$ws.Cells.Item(107, 14).Value = -4132.5

# diff hunk @42020 -> GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 16706.209
$ws.Cells.Item(132, 9).Value = 22341.182
$ws.Cells.Item(132, 10).Value = 11938.154
$ws.Cells.Item(132, 11).Value = 67023.546
$ws.Cells.Item(132, 12).Value = 35814.462
$ws.Cells.Item(132, 13).Value = -64493.546
$ws.Cells.Item(132, 14).Value = -40874.462

# diff hunk @42785 -> LTW!row6
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(6, 8).Value = 50000
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 50000
$ws.Cells.Item(6, 11).Value = 0
$ws.Cells.Item(6, 12).Value = 50000
$ws.Cells.Item(6, 14).Value = -50224

# diff hunk @43575 -> LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3048
$ws.Cells.Item(22, 9).Value = 2219.75
$ws.Cells.Item(22, 10).Value = 3462.125
$ws.Cells.Item(22, 11).Value = 2219.75
$ws.Cells.Item(22, 12).Value = 3462.125
$ws.Cells.Item(22, 13).Value = -1924.75
$ws.Cells.Item(22, 14).Value = -4052.125

# diff hunk @43814 -> LTW!row27
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 3048
$ws.Cells.Item(27, 9).Value = 2219.75
$ws.Cells.Item(27, 10).Value = 3462.125
$ws.Cells.Item(27, 11).Value = 2219.75
$ws.Cells.Item(27, 12).Value = 3462.125
$ws.Cells.Item(27, 13).Value = -2112.75
$ws.Cells.Item(27, 14).Value = -3676.125

# diff hunk @44730 -> LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 7142.7144
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 7142.7144
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 7142.7144
$ws.Cells.Item(46, 13).Value = $null
$ws.Cells.Item(46, 14).Value = -7518.7144

# diff hunk @45165 -> LTW!row55
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 2463.6875
$ws.Cells.Item(55, 9).Value = 1958.1111
$ws.Cells.Item(55, 10).Value = 3113.7144
$ws.Cells.Item(55, 11).Value = 1958.1111
$ws.Cells.Item(55, 12).Value = 3113.7144
$ws.Cells.Item(55, 13).Value = -1785.1111
$ws.Cells.Item(55, 14).Value = -3459.7144

# diff hunk @46976 -> LTW!row93
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 4674.1113
$ws.Cells.Item(93, 9).Value = 3017.25
$ws.Cells.Item(93, 10).Value = 5999.6
$ws.Cells.Item(93, 11).Value = 3017.25
$ws.Cells.Item(93, 12).Value = 5999.6
$ws.Cells.Item(93, 13).Value = -1769.25
$ws.Cells.Item(93, 14).Value = -8495.6

# diff hunk @48361 -> LTW!row122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 6594.2856
$ws.Cells.Item(122, 9).Value = 5632.2
$ws.Cells.Item(122, 10).Value = 8999.5
$ws.Cells.Item(122, 11).Value = 16896.6
$ws.Cells.Item(122, 12).Value = 26998.5
$ws.Cells.Item(122, 13).Value = -14446.6
$ws.Cells.Item(122, 14).Value = -31898.5

# diff hunk @48845 -> LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 8993818
$ws.Cells.Item(132, 9).Value = 16696950
$ws.Cells.Item(132, 10).Value = 6831
$ws.Cells.Item(132, 11).Value = 50090850
$ws.Cells.Item(132, 12).Value = 20493
$ws.Cells.Item(132, 13).Value = -50088320
$ws.Cells.Item(132, 14).Value = -25553

# diff hunk @52631 -> WVR!row69
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(69, 8).Value = 20271
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 20271
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 20271
$ws.Cells.Item(69, 14).Value = -21769

# diff hunk @52775 -> WVR!row72
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(72, 8).Value = 20271
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 20271
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).Value = 60813
$ws.Cells.Item(72, 14).Value = -68301

# diff hunk @52968 -> WVR!row76
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(76, 8).Value = 40000
$ws.Cells.Item(76, 9).Value = 40000
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 11).Value = 40000
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 13).Value = -39685

# diff hunk @53112 -> WVR!row79
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(79, 8).Value = 40000
$ws.Cells.Item(79, 9).Value = 40000
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 11).Value = 40000
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 13).Value = -38908

# diff hunk @54739 -> WVR!row113
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 5556469.5
$ws.Cells.Item(113, 9).Value = 722
$ws.Cells.Item(113, 10).Value = 15874287
$ws.Cells.Item(113, 11).Value = 2166
$ws.Cells.Item(113, 12).Value = 47622861
$ws.Cells.Item(113, 13).Value = 4
$ws.Cells.Item(113, 14).Value = -47627201

# diff hunk @55652 -> WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 23816372
$ws.Cells.Item(132, 9).Value = 27780766
$ws.Cells.Item(132, 10).Value = 30000
$ws.Cells.Item(132, 11).Value = 83342298
$ws.Cells.Item(132, 12).Value = 90000
$ws.Cells.Item(132, 13).Value = -83339768
$ws.Cells.Item(132, 14).Value = -95060
